$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("10per change")

# --- Fix E80: was stored as text "526371", should become a real number 526371 ---
$ws.Range("E80").Value = 526371

# --- Add row 81: PRESTIGE ---
$ws.Range("A81").Value = "27/06/2024 07:44:55"
$ws.Range("B81").Value = 1
$ws.Range("C81").Value = "PRESTIGE"
$ws.Range("D81").Value = "Prestige Estates Projects Limited"

# E81 must stay text (bsecode "533274" stored as a string, not a number)
$ws.Range("E81").Value = "'533274"
$ws.Range("E81").Style = "Normal"

$ws.Range("F81").Value = -3.22
$ws.Range("G81").Value = 1865.05
$ws.Range("H81").Value = 540053

# --- Add row 82: NMDC ---
$ws.Range("A82").Value = "27/06/2024 07:44:55"
$ws.Range("B82").Value = 2
$ws.Range("C82").Value = "NMDC"
$ws.Range("D82").Value = "Nmdc Limited"

# E82 must stay text (bsecode "526371" stored as a string, not a number)
$ws.Range("E82").Value = "'526371"
$ws.Range("E82").Style = "Normal"

$ws.Range("F82").Value = -2.35
$ws.Range("G82").Value = 243.4
$ws.Range("H82").Value = 9903555
